$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($addr, $val) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-TextCell "D2" '36.483.80'
Set-TextCell "E2" '  +1.71%  '
Set-TextCell "D3" '1.947.48'
Set-TextCell "E3" '  -0.28%  '
Set-TextCell "E4" '  +0.04%  '
Set-TextCell "D5" '243.48'
Set-TextCell "E5" '  +0.88%  '
Set-TextCell "D6" '0.616'
Set-TextCell "E6" '  -1.21%  '
Set-TextCell "D7" '58.22'
Set-TextCell "E7" '  -6.07%  '
Set-TextCell "E8" '  -0.04%  '
Set-TextCell "D9" '0.367'
Set-TextCell "E9" '  -0.08%  '
Set-TextCell "D10" '55.66'
Set-TextCell "E10" '  -0.60%  '
Set-TextCell "D11" '0.0837'
Set-TextCell "E11" '  +5.28%  '
Set-TextCell "D12" '0.103'
Set-TextCell "E12" '  +0.69%  '
Set-TextCell "D13" '21.61'
Set-TextCell "E13" '  -1.80%  '
Set-TextCell "D14" '0.824'
Set-TextCell "E14" '  -3.26%  '
Set-TextCell "D15" '2.235.03'
Set-TextCell "E15" '  -0.37%  '
Set-TextCell "D16" '13.58'
Set-TextCell "E16" '  -2.52%  '
Set-TextCell "D17" '5.24'
Set-TextCell "E17" '  -2.80%  '
Set-TextCell "D18" '1.942.44'
Set-TextCell "E18" '  -1.40%  '
Set-TextCell "D19" '36.417.01'
Set-TextCell "E19" '  +1.94%  '
Set-TextCell "D20" '69.82'
Set-TextCell "E20" '  -1.54%  '
Set-TextCell "D21" '0.0₃0863'
Set-TextCell "E21" '  +1.65%  '
Set-TextCell "D22" '229.61'
Set-TextCell "E22" '  -3.24%  '
Set-TextCell "D23" '5.07'
Set-TextCell "E23" '  -2.15%  '
Set-TextCell "E24" '  +0.08%  '
Set-TextCell "D25" '2.44'
Set-TextCell "E25" '  -2.95%  '
Set-TextCell "D26" '2.29'
Set-TextCell "E26" '  -0.05%  '
Set-TextCell "D27" '9.22'
Set-TextCell "E27" '  -5.78%  '
Set-TextCell "D28" '162.19'
Set-TextCell "E28" '  +1.96%  '
Set-TextCell "D29" '19.40'
Set-TextCell "E29" '  -1.68%  '
Set-TextCell "D30" '0.126'
Set-TextCell "E30" '  -0.94%  '
Set-TextCell "E31" '  -1.13%  '
Set-TextCell "D32" '1.15'
Set-TextCell "E32" '  +1.40%  '
Set-TextCell "D33" '4.68'
Set-TextCell "E33" '  -3.15%  '
Set-TextCell "D34" '0.0628'
Set-TextCell "E34" '  +1.87%  '
Set-TextCell "D35" '4.28'
Set-TextCell "E35" '  -2.26%  '
Set-TextCell "D36" '6.22'
Set-TextCell "E36" '  -0.35%  '
Set-TextCell "E37" '  -0.02%  '
Set-TextCell "D38" '1.78'
Set-TextCell "E38" '  -2.63%  '
Set-TextCell "D39" '2.14'
Set-TextCell "E39" '  -5.67%  '
Set-TextCell "D40" '3.04'
Set-TextCell "E40" '  -1.78%  '
Set-TextCell "E41" '  +0.46%  '
Set-TextCell "D42" '2.86'
Set-TextCell "E42" '  +2.39%  '
Set-TextCell "D43" '1.17'
Set-TextCell "E43" '  -3.60%  '
Set-TextCell "D44" '0.0209'
Set-TextCell "E44" '  -0.72%  '
Set-TextCell "D45" '16.09'
Set-TextCell "E45" '  +0.14%  '
Set-TextCell "D46" '1.351.04'
Set-TextCell "E46" '  +1.26%  '
Set-TextCell "D47" '1.03'
Set-TextCell "E47" '  -4.77%  '
Set-TextCell "D48" '87.78'
Set-TextCell "E48" '  -4.28%  '
Set-TextCell "D49" '7.18'
Set-TextCell "E49" '  -4.67%  '
Set-TextCell "E50" '  +2.11%  '
Set-TextCell "D51" '45.49'
Set-TextCell "E51" '  +4.09%  '
